$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to update for each data row (rows 2-25 correspond to bus indices 0-23)
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row, in the same order as $columns
$newValues = @{
    2 = @(1.02, 1.044851754993767, 1.04239141642672, 1.052316954061564, 1.062528347295586, 1.034161417093622, 1.049914625706421, 1.045168235852213, 1.055066021810953, 1.065249450043176, 1.051405624901815)
    3 = @(1.02, 1.046424936104267, 1.043517314444692, 1.05379078071919, 1.064259335791564, 1.034454904882206, 1.051132148410952, 1.046104130472947, 1.056350973924125, 1.066792970410825, 1.052624876628238)
    4 = @(1.02, 1.047440937957644, 1.044243957472834, 1.054742928570332, 1.065378205842134, 1.034642659653634, 1.051917664654829, 1.046707251274304, 1.057180400432468, 1.067790077958462, 1.053411508395281)
    5 = @(1.02, 1.047867607768563, 1.044548992052301, 1.055142856976782, 1.065848302978901, 1.034721078772074, 1.05224735194673, 1.046960218577341, 1.057528614744791, 1.068208875517705, 1.053741663880952)
    6 = @(1.02, 1.047939220971343, 1.044600182708736, 1.055209986171315, 1.065927218480187, 1.034734215665204, 1.052302676224252, 1.047002658776114, 1.057587053738716, 1.068279171079663, 1.053797066725294)
    7 = @(1.02, 1.047446640925568, 1.044248035105618, 1.054748273815916, 1.065384488374188, 1.034643709506884, 1.051922072079907, 1.04671063372633, 1.057185055152202, 1.067795675458907, 1.053415922079409)
    8 = @(1.02, 1.045383828926724, 1.042772313483339, 1.052815359834284, 1.063113596674762, 1.034261049392893, 1.050326574425108, 1.045485039716336, 1.055500701082468, 1.06577143796539, 1.051818158634934)
    9 = @(1.02, 1.041733502081688, 1.040157175113587, 1.049397320314196, 1.059102377303651, 1.033570182077322, 1.047497134803623, 1.043306251022439, 1.052516804873483, 1.062191389192016, 1.048984700884645)
    10 = @(1.02, 1.039288983563507, 1.038403471509772, 1.047110022071995, 1.056421113479119, 1.033098329629609, 1.045598288549796, 1.041840504685519, 1.050516389886313, 1.059795302107617, 1.047083158051337)
    11 = @(1.02, 1.038227739870717, 1.037641576209658, 1.046117435090232, 1.055258262926783, 1.032891308221909, 1.044772983924875, 1.041202604817874, 1.049647432452612, 1.058755411016825, 1.046256681399185)
    12 = @(1.02, 1.037833121348226, 1.037358186878398, 1.045748407341369, 1.054826038563261, 1.032814002028875, 1.044465955344516, 1.040965169358238, 1.049324237995713, 1.058368781691015, 1.045949216803025)
    13 = @(1.02, 1.037917787804078, 1.037418992506113, 1.045827580453892, 1.054918765478342, 1.032830603030842, 1.044531835601742, 1.041016122403889, 1.049393583690505, 1.058451731808094, 1.046015190617772)
    14 = @(1.02, 1.038195129327924, 1.037618159106839, 1.046086938071869, 1.055222541113121, 1.032884926431372, 1.044747614563478, 1.041182988386494, 1.049620725823412, 1.058723459679229, 1.046231276010385)
    15 = @(1.02, 1.038365951867077, 1.037740820654179, 1.046246691922864, 1.055409668514867, 1.032918342565308, 1.044880500057126, 1.041285734747455, 1.04976061898497, 1.058890831149988, 1.046364350216676)
    16 = @(1.02, 1.039359355824936, 1.038453982057839, 1.047175850153407, 1.056498248081751, 1.033112011716815, 1.045652995377856, 1.04188277152082, 1.050574000663298, 1.059864265353526, 1.047137942569368)
    17 = @(1.02, 1.039981746877138, 1.038900646807716, 1.047758097764285, 1.057180582292431, 1.033232768873863, 1.046136727083066, 1.042256409385232, 1.051083466580313, 1.060474232443267, 1.047622361229083)
    18 = @(1.02, 1.040344512431967, 1.039160935187492, 1.048097504115016, 1.057578399075255, 1.033302943553835, 1.046418581928496, 1.042474035625831, 1.051380363347589, 1.060829788078151, 1.047904616340724)
    19 = @(1.02, 1.040468161590744, 1.039249645699612, 1.048213197803448, 1.057714014638008, 1.033326827128186, 1.046514636933372, 1.042548188109642, 1.051481552630435, 1.060950985059499, 1.048000807754724)
    20 = @(1.02, 1.039914997665116, 1.038852749158579, 1.047695649781966, 1.057107392710881, 1.033219839781809, 1.04608485810554, 1.04221635373548, 1.051028833300834, 1.060408812397666, 1.04757041859166)
    21 = @(1.02, 1.03811347100113, 1.037559520263375, 1.046010573098162, 1.055133094864359, 1.032868940861925, 1.044684086188107, 1.041133864125266, 1.049553849935866, 1.058643452830537, 1.046167657417433)
    22 = @(1.02, 1.036978312616772, 1.03674417030146, 1.044949144925762, 1.053890091892611, 1.032645947802646, 1.043800620003811, 1.040450415396447, 1.048624005736854, 1.057531369357088, 1.045282936609821)
    23 = @(1.02, 1.037580319454782, 1.037176618085494, 1.045512016779206, 1.054549195106273, 1.032764386086881, 1.044269225461364, 1.040812996474627, 1.049117170501536, 1.058121111778688, 1.045752207540874)
    24 = @(1.02, 1.039945159578561, 1.038874392797407, 1.047723867982522, 1.057140464480055, 1.03322568268829, 1.046108296383704, 1.042234454118673, 1.051053520550441, 1.060438373604425, 1.047593890154867)
    25 = @(1.02, 1.042679086289162, 1.040835034880914, 1.050282442438823, 1.06014058013323, 1.03375076518789, 1.048230790495, 1.043871824103631, 1.053290142443086, 1.063118528652901, 1.049719398451255)
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
